$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$YELLOW = 65535

function Set-Val($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
}

function Set-Highlight($row, $col, [bool]$highlighted) {
    $c = $ws.Cells.Item($row, $col)
    if ($highlighted) {
        $c.Interior.Color = $YELLOW
    } else {
        $c.ClearFormats()
    }
}

# ---------- Row 22 (Fluid community 0) ----------
Set-Val 22 2 "[2, 21]"
Set-Val 22 3 0.413043
Set-Highlight 22 3 $true
Set-Val 22 4 "[13, 1, 9]"
Set-Val 22 5 0.347826
Set-Val 22 6 "[0, 0, 0, 1, 1, 0, 1, 0, 0, 1, 1, 1, 1, 1, 2, 0, 1, 1, 1, 2, 1, 1, 1, 1, 0, 1, 1, 1, 0, 1]"
Set-Val 22 7 0.355072
Set-Highlight 22 7 $false
Set-Val 22 8 "[1, 10, 4, 8]"
Set-Val 22 9 0.336957
Set-Highlight 22 9 $false

# ---------- Row 23 (Fluid community 1) ----------
Set-Val 23 2 "[12, 0]"
Set-Val 23 4 "[7, 0, 5]"
Set-Val 23 5 0.388889
Set-Highlight 23 5 $false
Set-Val 23 6 "[1, 2, 1, 0, 0, 1, 0, 1, 1, 0, 0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 1, 1, 0, 0, 0, 2, 0]"
Set-Val 23 7 0.711111
Set-Val 23 8 "[0, 3, 3, 6]"
Set-Val 23 9 0.375
Set-Highlight 23 9 $false

# ---------- Row 25 (Fluid community 3) ----------
Set-Val 25 2 "[2, 0]"
Set-Val 25 4 "[1, 0, 1]"
Set-Val 25 5 0.333333
Set-Highlight 25 5 $false
Set-Val 25 6 "[0, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0]"
Set-Val 25 7 0.9666670000000001
Set-Val 25 8 "[0, 0, 0, 2]"
Set-Val 25 9 0.75

# ---------- Row 26 (Fluid community 4) ----------
Set-Val 26 2 "[2, 0]"
Set-Val 26 3 0.5
Set-Highlight 26 3 $true
Set-Val 26 4 "[2, 0, 0]"
Set-Val 26 5 0.666667
Set-Highlight 26 5 $true
Set-Val 26 6 "[0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0]"
Set-Val 26 8 "[0, 2, 0, 0]"
